# Apply the "removed ER tags from non-ER templates and non-ER tags" edit.
#
# The metadata sheet ("SwateTemplateMetadata") carries an ER (endpoint
# repository) block in rows 8-10 (ER / ER Term Accession Number / ER Term
# Source REF). This template isn't tied to a specific ER, so the values
# are cleared, leaving the row labels (column A) and cell styles intact.
#
# The metadata sheet is also renamed to "isa_template" and becomes the
# active/selected sheet, matching the author's final workbook state.

$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item("SwateTemplateMetadata")

# Clear the ER list values (row 8: ER, row 9: ER Term Accession Number,
# row 10: ER Term Source REF) - keep the label column (A) untouched.
$metaSheet.Range("B8").ClearContents()
$metaSheet.Range("B9").ClearContents()
$metaSheet.Range("B10").ClearContents()

# This sheet becomes the active / selected tab in the saved workbook.
$metaSheet.Activate()

# Rename the metadata sheet.
$metaSheet.Name = "isa_template"
